$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing values in row 173 ---
$ws.Range("B173").Value = 117419
$ws.Range("F173").Value = 8604
$ws.Range("G173").Value = 5167
$ws.Range("H173").Value = 2743
$ws.Range("I173").Value = 1551
$ws.Range("R173").Value = 92277
$ws.Range("V173").Value = 7453
$ws.Range("W173").Value = 4587
$ws.Range("X173").Value = 2030
$ws.Range("Y173").Value = 1388

# --- Update existing values in row 174 ---
$ws.Range("B174").Value = 102419
$ws.Range("E174").Value = 36916
$ws.Range("F174").Value = 7722
$ws.Range("G174").Value = 4656
$ws.Range("I174").Value = 2454
$ws.Range("R174").Value = 81650
$ws.Range("U174").Value = 29067
$ws.Range("V174").Value = 6954
$ws.Range("W174").Value = 4006
$ws.Range("Y174").Value = 1895

# --- Update existing values in row 175 ---
$ws.Range("B175").Value = 94943
$ws.Range("D175").Value = 29765
$ws.Range("E175").Value = 29455
$ws.Range("F175").Value = 6894
$ws.Range("G175").Value = 5257
$ws.Range("H175").Value = 2356
$ws.Range("I175").Value = 2674
$ws.Range("R175").Value = 73469
$ws.Range("T175").Value = 22721
$ws.Range("U175").Value = 22846
$ws.Range("V175").Value = 5779
$ws.Range("W175").Value = 4441
$ws.Range("X175").Value = 2082
$ws.Range("Y175").Value = 2277

# --- Add new row 176 ---
# A176 must hold the literal text "01-07-2021" as a shared string, matching
# the style of the existing date-label cells in column A (no explicit cell
# style). Assigning that text straight to .Value/.Formula triggers Excel's
# locale-aware "looks like a date" auto-conversion (turning it into a date
# serial + a new number format). Routing it through a scratch formula cell
# and a values-only paste avoids that smart re-typing and lands the literal
# text, unformatted, exactly like the other rows.
$ws.Range("AA1").Formula = "=""01-07-2021"""
$ws.Range("AA1").Copy()
$ws.Range("A176").PasteSpecial(-4163)
$ws.Range("AA1").ClearContents()

$ws.Range("B176").Value = 88078
$ws.Range("C176").Value = 17689
$ws.Range("D176").Value = 22101
$ws.Range("E176").Value = 30892
$ws.Range("F176").Value = 6904
$ws.Range("G176").Value = 4865
$ws.Range("H176").Value = 3511
$ws.Range("I176").Value = 2116
$ws.Range("J176").Value = 18554
$ws.Range("K176").Value = 4426
$ws.Range("L176").Value = 5691
$ws.Range("M176").Value = 6202
$ws.Range("N176").Value = 666
$ws.Range("O176").Value = 517
$ws.Range("P176").Value = 901
$ws.Range("Q176").Value = 150
$ws.Range("R176").Value = 69524
$ws.Range("S176").Value = 13263
$ws.Range("T176").Value = 16410
$ws.Range("U176").Value = 24690
$ws.Range("V176").Value = 6237
$ws.Range("W176").Value = 4348
$ws.Range("X176").Value = 2610
$ws.Range("Y176").Value = 1966
